$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Delete the empty spacer column Q. Column Q had no data in it (it was a
#    blank gap between column P and the "*_edit"/"*_update" helper columns
#    that used to start at R); deleting it shifts everything from R onward
#    one column to the left (R->Q, S->R, ... AI->AH).
# ---------------------------------------------------------------------------
$ws.Range("Q:Q").Delete()

# ---------------------------------------------------------------------------
# 2. Keep the conditional-formatting target ranges in sync with the shift.
#    (A plain column delete does not renumber the sqref of existing
#    conditional-formatting rules, so do it explicitly.)
# ---------------------------------------------------------------------------

# -- Rules that used to target S1:S1048576 / W1:W1048576 / Y1 / AA1 /
#    AC1:AC1048576 / AH1:AH1048576 / Y4:Y1048576 / AA4:AA1048576 /
#    U1:U1048576 now live one column to the left. The column delete left the
#    rule objects anchored to their old (pre-shift) addresses, so fetch them
#    from S1 (where they still live) and re-point them at R1.
$blk1 = $ws.Range("S1").FormatConditions
$blk1Rule1 = $blk1.Item(1)
$blk1Rule2 = $blk1.Item(2)
$blk1Rule1.ModifyAppliesToRange($ws.Range("R1:R1048576"))
$blk1Rule2.ModifyAppliesToRange($ws.Range("R1:R1048576"))
$extra1 = @("V1:V1048576", "X1", "Z1", "AB1:AB1048576", "AG1:AG1048576", "X4:X1048576", "Z4:Z1048576", "T1:T1048576")
foreach ($a in $extra1) {
    $rng = $ws.Range($a)
    $r1 = $rng.FormatConditions.Add(1, 3, '="updated"')
    $r1.Font.Color = 12611584
    $r1.Interior.Color = 16641439
    $r2 = $rng.FormatConditions.Add(1, 3, '="filled in"')
    $r2.Font.Color = 24832
    $r2.Interior.Color = 13561798
}

# -- Rule that used to target Y2:Y3 now lives at X2:X3.
$blk2 = $ws.Range("Y2").FormatConditions
$blk2.Item(1).ModifyAppliesToRange($ws.Range("X2:X3"))
$blk2.Item(2).ModifyAppliesToRange($ws.Range("X2:X3"))

# -- Rule that used to target AA2:AA3 now lives at Z2:Z3.
$blk3 = $ws.Range("AA2").FormatConditions
$blk3.Item(1).ModifyAppliesToRange($ws.Range("Z2:Z3"))
$blk3.Item(2).ModifyAppliesToRange($ws.Range("Z2:Z3"))

# ---------------------------------------------------------------------------
# 3. The cell that used to be AH5 (now AG5 after the shift) held the literal
#    string "unchanged". Replace it with an (unresolvable) array formula
#    referencing the name "unchanged" -- Excel can't resolve it, so it
#    evaluates to a #NAME? error.
# ---------------------------------------------------------------------------
$ws.Range("AG5").FormulaArray = "=unchanged"
